$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart | Roof Tile
$ws.Range("H19").Value = 4751.909
$ws.Range("I19").Value = 245.5
$ws.Range("J19").Value = 5753.3335
$ws.Range("K19").Value = 245.5
$ws.Range("L19").Value = 5753.3335
$ws.Range("M19").Value = -70.5
$ws.Range("N19").Value = -6103.3335

# Row 32: Automata for the People | Crab Oil
$ws.Range("H32").Value = 1850.5555
$ws.Range("I32").Value = 1500.5
$ws.Range("J32").Value = 1950.5714
$ws.Range("K32").Value = 1500.5
$ws.Range("L32").Value = 1950.5714
$ws.Range("M32").Value = -1174.5
$ws.Range("N32").Value = -2602.5714

# Row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 2932.8572
$ws.Range("J98").Value = 6999.6665
$ws.Range("L98").Value = 6999.6665
$ws.Range("N98").Value = -9995.666499999999

# Row 116: Growing Up | Growth Formula Kappa
$ws.Range("H116").Value = 3153.4814
$ws.Range("I116").Value = 2663.9333
$ws.Range("J116").Value = 3765.4167
$ws.Range("K116").Value = 2663.9333
$ws.Range("L116").Value = 3765.4167
$ws.Range("M116").Value = 778.0666999999999
$ws.Range("N116").Value = -10649.4167

# Row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 2932.8572
$ws.Range("J122").Value = 6999.6665
$ws.Range("L122").Value = 20998.9995
$ws.Range("N122").Value = -25898.9995

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3713.1135
$ws.Range("I138").Value = 1717.2222
$ws.Range("J138").Value = 6883.0586
$ws.Range("K138").Value = 5151.6666
$ws.Range("L138").Value = 20649.1758
$ws.Range("M138").Value = -11.66659999999956
$ws.Range("N138").Value = -30929.1758

$ws = $wb.Worksheets.Item("ARM")
# Row 26: Night Squawker | Iron Lantern Shield
$ws.Range("H26").Value = 606.7143
$ws.Range("I26").Value = 606.7143
$ws.Range("K26").Value = 606.7143
$ws.Range("M26").Value = -276.7143

# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 5882.636
$ws.Range("I32").Value = 5018.033
$ws.Range("J32").Value = 7735.357
$ws.Range("K32").Value = 5018.033
$ws.Range("L32").Value = 7735.357
$ws.Range("M32").Value = -4731.033
$ws.Range("N32").Value = -8309.357

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 4120.0435
$ws.Range("I61").Value = 643.625
$ws.Range("J61").Value = 5974.1333
$ws.Range("K61").Value = 643.625
$ws.Range("L61").Value = 5974.1333
$ws.Range("M61").Value = -431.625
$ws.Range("N61").Value = -6398.1333

# Row 102: Smells of Rich Tama-hagane | Tama-hagane Ingot
$ws.Range("H102").Value = 2129.5806
$ws.Range("I102").Value = 1782.1923
$ws.Range("K102").Value = 1782.1923
$ws.Range("M102").Value = -160.1922999999999

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 4120.0435
$ws.Range("I136").Value = 643.625
$ws.Range("J136").Value = 5974.1333
$ws.Range("K136").Value = 1930.875
$ws.Range("L136").Value = 17922.3999
$ws.Range("M136").Value = 619.125
$ws.Range("N136").Value = -23022.3999

$ws = $wb.Worksheets.Item("BSM")
# Row 21: Awl or Nothing | Iron Awl
$ws.Range("H21").Value = 12266
$ws.Range("J21").Value = 12266
$ws.Range("L21").Value = 12266
$ws.Range("N21").Value = -12738

# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 5048.6665
$ws.Range("I86").Value = 4373
$ws.Range("J86").Value = 6400
$ws.Range("K86").Value = 4373
$ws.Range("L86").Value = 6400
$ws.Range("M86").Value = -3250
$ws.Range("N86").Value = -8646

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 5048.6665
$ws.Range("I89").Value = 4373
$ws.Range("J89").Value = 6400
$ws.Range("K89").Value = 21865
$ws.Range("L89").Value = 32000
$ws.Range("M89").Value = -16249
$ws.Range("N89").Value = -43232

$ws = $wb.Worksheets.Item("CRP")
# Row 6: Got Your Back | Square Maple Shield
$ws.Range("H6").Value = 606794.1
$ws.Range("I6").Value = 857705.9
$ws.Range("J6").Value = 21333.334
$ws.Range("K6").Value = 857705.9
$ws.Range("L6").Value = 21333.334
$ws.Range("M6").Value = -857592.9
$ws.Range("N6").Value = -21559.334

# Row 10: Spears and Sorcery | Maple Crook
$ws.Range("H10").Value = 23009.857
$ws.Range("I10").Value = 351
$ws.Range("J10").Value = 40004
$ws.Range("K10").Value = 351
$ws.Range("L10").Value = 40004
$ws.Range("M10").Value = -212
$ws.Range("N10").Value = -40282

# Row 19: Shielding Sales | Square Ash Shield
$ws.Range("H19").Value = 948.7273
$ws.Range("I19").Value = 266
$ws.Range("J19").Value = 1517.6666
$ws.Range("K19").Value = 266
$ws.Range("L19").Value = 1517.6666
$ws.Range("M19").Value = -96
$ws.Range("N19").Value = -1857.6666

# Row 24: What You Need | Square Ash Shield
$ws.Range("H24").Value = 948.7273
$ws.Range("I24").Value = 266
$ws.Range("J24").Value = 1517.6666
$ws.Range("K24").Value = 266
$ws.Range("L24").Value = 1517.6666
$ws.Range("M24").Value = -96
$ws.Range("N24").Value = -1857.6666

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 3084.7778
$ws.Range("J31").Value = 3314
$ws.Range("L31").Value = 3314
$ws.Range("N31").Value = -3904

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 3084.7778
$ws.Range("J34").Value = 3314
$ws.Range("L34").Value = 3314
$ws.Range("N34").Value = -3718

# Row 99: O Pine | Pine Lumber
$ws.Range("H99").Value = 1867.32
$ws.Range("I99").Value = 1200
$ws.Range("J99").Value = 2078.0527
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 2078.0527
$ws.Range("M99").Value = 298
$ws.Range("N99").Value = -5074.0527

# Row 105: Zelkova, My Love | Zelkova Lumber
$ws.Range("H105").Value = 3640.4666
$ws.Range("I105").Value = 3500.75
$ws.Range("J105").Value = 4199.3335
$ws.Range("K105").Value = 3500.75
$ws.Range("L105").Value = 4199.3335
$ws.Range("M105").Value = -1753.75
$ws.Range("N105").Value = -7693.3335

# Row 126: A Better Conductor | Red Pine Lumber
$ws.Range("H126").Value = 1867.32
$ws.Range("I126").Value = 1200
$ws.Range("J126").Value = 2078.0527
$ws.Range("K126").Value = 3600
$ws.Range("L126").Value = 6234.158100000001
$ws.Range("M126").Value = -1130
$ws.Range("N126").Value = -11174.1581

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 3375.6296
$ws.Range("I132").Value = 2084
$ws.Range("J132").Value = 5254.364
$ws.Range("K132").Value = 6252
$ws.Range("L132").Value = 15763.092
$ws.Range("M132").Value = -3722
$ws.Range("N132").Value = -20823.092

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 2448.12
$ws.Range("I134").Value = 1440.7059
$ws.Range("J134").Value = 4588.875
$ws.Range("K134").Value = 4322.1177
$ws.Range("L134").Value = 13766.625
$ws.Range("M134").Value = -1787.1177
$ws.Range("N134").Value = -18836.625

$ws = $wb.Worksheets.Item("CUL")
# Row 10: A Real Fungi | Chanterelle Saute
$ws.Range("H10").Value = 2405.889
$ws.Range("I10").Value = 216.66667
$ws.Range("J10").Value = 3500.5
$ws.Range("K10").Value = 650.00001
$ws.Range("L10").Value = 10501.5
$ws.Range("M10").Value = -511.00001
$ws.Range("N10").Value = -10779.5

# Row 29: For Crumbs' Sake | Honey Muffin
$ws.Range("H29").Value = 26650
$ws.Range("I29").Value = 100
$ws.Range("J29").Value = 35500
$ws.Range("K29").Value = 300
$ws.Range("L29").Value = 106500
$ws.Range("M29").Value = -23
$ws.Range("N29").Value = -107054

# Row 47: Winter of Our Discontent | Mugwort Carp
$ws.Range("H47").Value = 1946
$ws.Range("I47").Value = 219
$ws.Range("J47").Value = 3097.3333
$ws.Range("K47").Value = 657
$ws.Range("L47").Value = 9291.999899999999
$ws.Range("M47").Value = -226
$ws.Range("N47").Value = -10153.9999

# Row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 2613
$ws.Range("I132").Value = 1062.8334
$ws.Range("J132").Value = 5713.3335
$ws.Range("K132").Value = 9565.500599999999
$ws.Range("L132").Value = 51420.0015
$ws.Range("M132").Value = -7035.500599999999
$ws.Range("N132").Value = -56480.0015

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle | Durium Ingot
$ws.Range("H102").Value = 1721.7241
$ws.Range("I102").Value = 1018.6957
$ws.Range("J102").Value = 4416.6665
$ws.Range("K102").Value = 1018.6957
$ws.Range("L102").Value = 4416.6665
$ws.Range("M102").Value = 603.3043
$ws.Range("N102").Value = -7660.6665

# Row 125: Pewter-hewn Punishment | Pewter Choker of Slaying
$ws.Range("H125").Value = 23260
$ws.Range("J125").Value = 23260
$ws.Range("L125").Value = 23260
$ws.Range("N125").Value = -28180

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban | Leather
$ws.Range("H7").Value = 2024.0588
$ws.Range("I7").Value = 1401.909
$ws.Range("J7").Value = 3164.6667
$ws.Range("K7").Value = 1401.909
$ws.Range("L7").Value = 3164.6667
$ws.Range("M7").Value = -1289.909
$ws.Range("N7").Value = -3388.6667

# Row 40: Best Served Toad | Toad Leather
$ws.Range("H40").Value = 2534.7778
$ws.Range("I40").Value = 1154
$ws.Range("J40").Value = 2929.2856
$ws.Range("K40").Value = 1154
$ws.Range("L40").Value = 2929.2856
$ws.Range("M40").Value = -1018
$ws.Range("N40").Value = -3201.2856

# Row 106: If the Shoe Fits | Gazelleskin Boots of Casting
$ws.Range("H106").Value = 17674
$ws.Range("J106").Value = 17674
$ws.Range("L106").Value = 17674
$ws.Range("N106").Value = -20198

# Row 126: Battered Books | Saiga Leather
$ws.Range("H126").Value = 2024.0588
$ws.Range("I126").Value = 1401.909
$ws.Range("J126").Value = 3164.6667
$ws.Range("K126").Value = 4205.727000000001
$ws.Range("L126").Value = 9494.000100000001
$ws.Range("M126").Value = -1735.727000000001
$ws.Range("N126").Value = -14434.0001

# Row 127: Loyal Turncoat | Saigaskin Coat of Fending
$ws.Range("H127").Value = 28894
$ws.Range("J127").Value = 28894
$ws.Range("L127").Value = 28894
$ws.Range("N127").Value = -38814

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 6457.5
$ws.Range("I132").Value = 1940.55
$ws.Range("K132").Value = 5821.65
$ws.Range("M132").Value = -3291.65
